# Update cryptocurrency price/volume data on the active worksheet
# to reflect the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal text. Numeric-looking strings
# (e.g. "213.61") are given a leading apostrophe so Excel keeps them
# as text instead of silently converting them to a Double, matching
# the original inline-string cell content exactly.
function Set-TextValue($cell, $text) {
    if ($text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

Set-TextValue $ws.Range("D2") "27.915.75"
Set-TextValue $ws.Range("D3") "1.644.20"
Set-TextValue $ws.Range("E3") "  +1.19%  "
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "213.61"
Set-TextValue $ws.Range("E5") "  +0.97%  "
Set-TextValue $ws.Range("E6") "  -0.40%  "
Set-TextValue $ws.Range("E7") "  +0.02%  "
Set-TextValue $ws.Range("E8") "  +1.49%  "
Set-TextValue $ws.Range("E9") "  +0.11%  "
Set-TextValue $ws.Range("E10") "  +0.79%  "
Set-TextValue $ws.Range("D11") "0.0876"
Set-TextValue $ws.Range("E11") "  -1.60%  "
Set-TextValue $ws.Range("D12") "1.877.30"
Set-TextValue $ws.Range("E12") "  +1.21%  "
Set-TextValue $ws.Range("D13") "1.645.79"
Set-TextValue $ws.Range("E13") "  +1.31%  "
Set-TextValue $ws.Range("D14") "0.575"
Set-TextValue $ws.Range("E14") "  +4.49%  "
Set-TextValue $ws.Range("E15") "  +0.43%  "
Set-TextValue $ws.Range("D16") "65.85"
Set-TextValue $ws.Range("E16") "  +0.97%  "
Set-TextValue $ws.Range("D17") "27.898.11"
Set-TextValue $ws.Range("D18") "230.15"
Set-TextValue $ws.Range("E18") "  -0.84%  "
Set-TextValue $ws.Range("E19") "  +0.80%  "
Set-TextValue $ws.Range("D20") "7.63"
Set-TextValue $ws.Range("E20") "  +0.84%  "
Set-TextValue $ws.Range("E21") "  +0.02%  "
Set-TextValue $ws.Range("E22") "  +4.35%  "
Set-TextValue $ws.Range("E23") "  +1.30%  "
Set-TextValue $ws.Range("E24") "  +1.97%  "
Set-TextValue $ws.Range("E26") "  +0.61%  "
Set-TextValue $ws.Range("E27") "  +0.93%  "
Set-TextValue $ws.Range("D28") "15.71"
Set-TextValue $ws.Range("E28") "  +0.89%  "
Set-TextValue $ws.Range("E29") "  +0.04%  "
Set-TextValue $ws.Range("E30") "  +1.04%  "
Set-TextValue $ws.Range("E31") "  +0.17%  "
Set-TextValue $ws.Range("E32") "  +1.89%  "
Set-TextValue $ws.Range("D33") "1.424.70"
Set-TextValue $ws.Range("E33") "  -2.92%  "
Set-TextValue $ws.Range("D34") "3.10"
Set-TextValue $ws.Range("E34") "  +0.76%  "
Set-TextValue $ws.Range("E36") "  -0.20%  "
Set-TextValue $ws.Range("E37") "  +1.51%  "
Set-TextValue $ws.Range("E38") "  +0.60%  "
Set-TextValue $ws.Range("D39") "0.926"
Set-TextValue $ws.Range("E39") "  -2.65%  "
Set-TextValue $ws.Range("D40") "0.558"
Set-TextValue $ws.Range("E40") "  +0.22%  "
Set-TextValue $ws.Range("D41") "1.03"
Set-TextValue $ws.Range("E41") "  +2.23%  "
Set-TextValue $ws.Range("E42") "  +0.00%  "
Set-TextValue $ws.Range("D43") "68.55"
Set-TextValue $ws.Range("E43") "  +1.13%  "
Set-TextValue $ws.Range("E44") "  +0.68%  "
Set-TextValue $ws.Range("E45") "  +2.75%  "
Set-TextValue $ws.Range("E46") "  +2.66%  "
Set-TextValue $ws.Range("E47") "  +0.06%  "
Set-TextValue $ws.Range("D48") "1.785.92"
Set-TextValue $ws.Range("E48") "  +1.19%  "
Set-TextValue $ws.Range("D49") "88.94"
Set-TextValue $ws.Range("E49") "  +1.70%  "
Set-TextValue $ws.Range("E50") "  +0.00%  "
Set-TextValue $ws.Range("B51") "Cronos"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.0506"
Set-TextValue $ws.Range("E51") "  +0.57%  "
